$p = $ppt.ActivePresentation
$s = $p.Slides.Add(6, 12)
$shp = $s.Shapes.AddShape(1, 72, 42, 24, 120)
$shp.Fill.Pattern = 26
$shp.Fill.BackColor.RGB = 413090
$shp.Fill.ForeColor.RGB = 1
Write-Output "Type after both colors set (back then fore):"
Write-Output $shp.Fill.Type
$shp.Fill.Pattern = 26
Write-Output "done"
